$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update "Last Update" date in A2 (43521 -> 43530) ---
$ws.Range("A2").Value = 43530

# --- Fill in meeting rows 44-47 (previously blank placeholder rows) ---
$ws.Range("A44").Value = "Team"
$ws.Range("B44").Value = "Wednesday, February, 27, 2019"
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = 1
$ws.Range("E44").Value = 1

$ws.Range("A45").Value = "Team"
$ws.Range("B45").Value = "Monday, March 4, 2019"
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 1

$ws.Range("A46").Value = "Team"
$ws.Range("B46").Value = "Wednesday, March 6, 2019"
$ws.Range("C46").Value = 1
$ws.Range("D46").Value = 1
$ws.Range("E46").Value = 1

$ws.Range("A47").Value = "Faculty Adv"
$ws.Range("B47").Value = "Wednesday, March 6, 2019"
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = 1

# --- Insert two blank rows before the "Ad hoc Meetings" section (old row 50) ---
$ws.Rows("50:51").Insert()

# --- Remove the now-surplus trailing blank rows in the ad hoc section ---
# After the insert above, the footer/total rows and the extra blank rows
# shifted down by two; trim five of the blank rows so the sheet ends at
# row 63 again, keeping the one row that carries the stray "s=15" style.
$ws.Rows(64).Delete()
$ws.Rows(63).Delete()
$ws.Rows(62).Delete()
$ws.Rows(60).Delete()
$ws.Rows(59).Delete()

# --- Populate the newly surfaced ad hoc meeting row 58 ---
$ws.Range("A58").Value = "Audio Pass through"
$ws.Range("B58").Value = "Sunday, March 10, 2019"

# --- Fix the leftover stray number-format style on C59 so it matches D59/E59 ---
$ws.Range("D59").Copy()
$ws.Range("C59").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update sheet selection to match the saved view ---
$ws.Range("F34").Select()
